$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Insert new header cell "id" in A1 of Sheet2 (shifts nothing else, just adds A1)
$ws2.Range("A1").Value = "id"

# Update the active selection on Sheet2 to A2
$ws2.Activate()
$ws2.Range("A2").Select()
